$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for every data row.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 302) {
    $lastRow = 302
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
